$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "XLookUp": add a VLOOKUP helper column in B (B3:B6) that
# looks up the Email (col P) from the EmployeeID table in H:P.
# ---------------------------------------------------------------
$wsX = $wb.Worksheets.Item("XLookUp")

$wsX.Range("B3").Formula = '=VLOOKUP(A3, $H$2:$P$10,9, FALSE)'
$wsX.Range("B4:B6").Formula = '=VLOOKUP(A4, $H$2:$P$10,9, FALSE)'

# Widen column A (new bestfit-ish width) and column B (now holds long emails)
$wsX.Columns.Item(1).ColumnWidth = 15.333333333333334
$wsX.Columns.Item(2).ColumnWidth = 43.666666666666664

$wsX.Activate()
$wsX.Range("B12").Select()

# ---------------------------------------------------------------
# Sheet "VLookUp": remove the unused/blank "Address" column (K),
# which shifts JobTitle..Email left from L:P into K:O, then add the
# analogous VLOOKUP helper column in B (B3:B5), plus a one-line note
# about FALSE/exact-match lookups in A12.
# ---------------------------------------------------------------
$wsV = $wb.Worksheets.Item("VLookUp")

$wsV.Columns("K:K").Delete()

$wsV.Range("B3").Formula = '=VLOOKUP(A3, H2:O10,8, FALSE)'
$wsV.Range("B4:B5").Formula = '=VLOOKUP(A4, H3:O11,8, FALSE)'

$wsV.Range("A12").Value = "Set to false, lookup value must be the first column of the table selected"

# Widen the JobTitle column (now H) and the Email column (now O)
$wsV.Columns.Item(8).ColumnWidth = 15.333333333333334
$wsV.Columns.Item(15).ColumnWidth = 43.666666666666664

$wsV.Activate()
$wsV.Range("A13").Select()
